$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 808.94446
$ws.Range("I15").Value = 808.94446
$ws.Range("K15").Value = 2426.83338
$ws.Range("M15").Value = -2257.83338
$ws.Range("H28").Value = 597.375
$ws.Range("I28").Value = 404.57144
$ws.Range("K28").Value = 404.57144
$ws.Range("M28").Value = 80.42856
$ws.Range("H33").Value = 257.55
$ws.Range("J33").Value = 540
$ws.Range("L33").Value = 540
$ws.Range("N33").Value = -998
$ws.Range("H62").Value = 10687.895
$ws.Range("J62").Value = 8867.166999999999
$ws.Range("L62").Value = 8867.166999999999
$ws.Range("N62").Value = -10115.167
$ws.Range("H65").Value = 10687.895
$ws.Range("J65").Value = 8867.166999999999
$ws.Range("L65").Value = 44335.835
$ws.Range("N65").Value = -50575.835
$ws.Range("H94").Value = 9600.714
$ws.Range("I94").Value = 2867.5
$ws.Range("K94").Value = 2867.5
$ws.Range("M94").Value = -2416.5
$ws.Range("I111").Value = 4137.5
$ws.Range("J111").Value = 7971.4287
$ws.Range("K111").Value = 12412.5
$ws.Range("L111").Value = 23914.2861
$ws.Range("M111").Value = -9345.5
$ws.Range("N111").Value = -30048.2861
$ws.Range("H129").Value = 2833.8333
$ws.Range("I129").Value = 2910
$ws.Range("J129").Value = 2714.1428
$ws.Range("K129").Value = 8730
$ws.Range("L129").Value = 8142.428400000001
$ws.Range("M129").Value = -3730
$ws.Range("N129").Value = -18142.4284
$ws.Range("H137").Value = 250004750
$ws.Range("I137").Value = 500000500
$ws.Range("K137").Value = 1500001500
$ws.Range("M137").Value = -1499998950
$ws.Range("H138").Value = 3307.5117
$ws.Range("I138").Value = 2196
$ws.Range("J138").Value = 3487.7568
$ws.Range("K138").Value = 6588
$ws.Range("L138").Value = 10463.2704
$ws.Range("M138").Value = -1448
$ws.Range("N138").Value = -20743.2704
$ws.Range("H141").Value = 1000
$ws.Range("I141").Value = 1000
$ws.Range("K141").Value = 3000
$ws.Range("M141").Value = 2180

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1452.579
$ws.Range("I32").Value = 1486.1351
$ws.Range("K32").Value = 1486.1351
$ws.Range("M32").Value = -1199.1351
$ws.Range("H45").Value = 2719.6
$ws.Range("I45").Value = 4550
$ws.Range("K45").Value = 4550
$ws.Range("M45").Value = -4173
$ws.Range("H61").Value = 3554.6667
$ws.Range("I61").Value = 2709.6667
$ws.Range("K61").Value = 2709.6667
$ws.Range("M61").Value = -2497.6667
$ws.Range("H102").Value = 3284.3635
$ws.Range("I102").Value = 3389.8333
$ws.Range("J102").Value = 3157.8
$ws.Range("K102").Value = 3389.8333
$ws.Range("L102").Value = 3157.8
$ws.Range("M102").Value = -1767.8333
$ws.Range("N102").Value = -6401.8
$ws.Range("H122").Value = 1751.5
$ws.Range("I122").Value = 1546.2778
$ws.Range("K122").Value = 4638.8334
$ws.Range("M122").Value = -2188.8334
$ws.Range("H132").Value = 2548.697
$ws.Range("I132").Value = 2266.6072
$ws.Range("K132").Value = 6799.821599999999
$ws.Range("M132").Value = -4269.821599999999
$ws.Range("H136").Value = 3554.6667
$ws.Range("I136").Value = 2709.6667
$ws.Range("K136").Value = 8129.000100000001
$ws.Range("M136").Value = -5579.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2222.6365
$ws.Range("I20").Value = 1999.75
$ws.Range("J20").Value = 2350
$ws.Range("K20").Value = 1999.75
$ws.Range("L20").Value = 2350
$ws.Range("M20").Value = -1752.75
$ws.Range("N20").Value = -2844
$ws.Range("H134").Value = 2239.468
$ws.Range("I134").Value = 1316
$ws.Range("K134").Value = 3948
$ws.Range("M134").Value = -1413

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4910.6924
$ws.Range("I58").Value = 2187.7144
$ws.Range("K58").Value = 2187.7144
$ws.Range("M58").Value = -1984.7144
$ws.Range("H107").Value = 708.8
$ws.Range("I107").Value = 619.4286
$ws.Range("J107").Value = 787
$ws.Range("K107").Value = 619.4286
$ws.Range("L107").Value = 787
$ws.Range("M107").Value = 1300.5714
$ws.Range("N107").Value = -4627
$ws.Range("H122").Value = 1891.6428
$ws.Range("I122").Value = 1462.091
$ws.Range("K122").Value = 4386.272999999999
$ws.Range("M122").Value = -1936.272999999999
$ws.Range("H132").Value = 86961920
$ws.Range("I132").Value = 142861200
$ws.Range("K132").Value = 428583600
$ws.Range("M132").Value = -428581070
$ws.Range("H136").Value = 4910.6924
$ws.Range("I136").Value = 2187.7144
$ws.Range("K136").Value = 6563.1432
$ws.Range("M136").Value = -4013.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 3499.5
$ws.Range("I112").Value = 3499.5
$ws.Range("K112").Value = 10498.5
$ws.Range("M112").Value = -9390.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8979.6
$ws.Range("I70").Value = 8132.6665
$ws.Range("K70").Value = 8132.6665
$ws.Range("M70").Value = -7862.6665
$ws.Range("H73").Value = 8979.6
$ws.Range("I73").Value = 8132.6665
$ws.Range("K73").Value = 8132.6665
$ws.Range("M73").Value = -7196.6665
$ws.Range("H97").Value = 1604.8695
$ws.Range("I97").Value = 1525.7368
$ws.Range("K97").Value = 1525.7368
$ws.Range("M97").Value = -1029.7368
$ws.Range("H102").Value = 2484.3333
$ws.Range("I102").Value = 2484.3333
$ws.Range("K102").Value = 2484.3333
$ws.Range("M102").Value = -862.3332999999998
$ws.Range("H104").Value = 80447.164
$ws.Range("J104").Value = 100670.75
$ws.Range("L104").Value = 100670.75
$ws.Range("N104").Value = -107658.75
$ws.Range("H113").Value = 3758.111
$ws.Range("I113").Value = 1266.6666
$ws.Range("J113").Value = 5003.8335
$ws.Range("K113").Value = 1266.6666
$ws.Range("L113").Value = 5003.8335
$ws.Range("M113").Value = 903.3334
$ws.Range("N113").Value = -9343.833500000001
$ws.Range("H132").Value = 17872468
$ws.Range("I132").Value = 22742210
$ws.Range("J132").Value = 16751.084
$ws.Range("K132").Value = 68226630
$ws.Range("L132").Value = 50253.25199999999
$ws.Range("M132").Value = -68224100
$ws.Range("N132").Value = -55313.25199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1164.6111
$ws.Range("I61").Value = 857.25
$ws.Range("J61").Value = 1779.3334
$ws.Range("K61").Value = 857.25
$ws.Range("L61").Value = 1779.3334
$ws.Range("M61").Value = -655.25
$ws.Range("N61").Value = -2183.3334
$ws.Range("H93").Value = 1682.5555
$ws.Range("I93").Value = 1682.5555
$ws.Range("K93").Value = 1682.5555
$ws.Range("M93").Value = -434.5554999999999
$ws.Range("H113").Value = 1164.6111
$ws.Range("I113").Value = 857.25
$ws.Range("J113").Value = 1779.3334
$ws.Range("K113").Value = 857.25
$ws.Range("L113").Value = 1779.3334
$ws.Range("M113").Value = 1312.75
$ws.Range("N113").Value = -6119.3334
$ws.Range("H122").Value = 4750.4375
$ws.Range("I122").Value = 4000.5833
$ws.Range("K122").Value = 12001.7499
$ws.Range("M122").Value = -9551.749899999999
$ws.Range("H132").Value = 2701.2917
$ws.Range("I132").Value = 2468.3333
$ws.Range("K132").Value = 7404.999899999999
$ws.Range("M132").Value = -4874.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8566.333000000001
$ws.Range("I81").Value = 10079.6
$ws.Range("K81").Value = 20159.2
$ws.Range("M81").Value = -19098.2
$ws.Range("H84").Value = 8566.333000000001
$ws.Range("I84").Value = 10079.6
$ws.Range("K84").Value = 100796
$ws.Range("M84").Value = -95492
$ws.Range("H136").Value = 6062911
$ws.Range("I136").Value = 6175002
$ws.Range("J136").Value = 9997
$ws.Range("K136").Value = 18525006
$ws.Range("L136").Value = 29991
$ws.Range("M136").Value = -18522456
$ws.Range("N136").Value = -35091
